$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the existing table (contents + formatting) so the sheet starts fresh.
$ws.UsedRange.Clear()

# Cells that carry the bold/bordered/centered "header" style: column headers
# in row 1 plus the index column A (two contiguous ranges - a multi-area
# Range only applies formatting to its first area via this COM shim).
$headerRange = $ws.Range("B1:E1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$indexRange = $ws.Range("A2:A6")
$indexRange.Font.Bold = $true
$indexRange.Borders.LineStyle = 1
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160

# New column headers (row 1, columns B:E).
$headers = @("mean:var", "mean:iqr", "mean:rvar", "mean:skew")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Row index values (A2:A6).
for ($r = 0; $r -lt 5; $r++) {
    $ws.Cells.Item($r + 2, 1).Value = $r
}

# Data values (B2:E6).
$data = @(
    @("-0.28**", "-0.42***", "-0.48***", "-0.02"),
    @("-0.42***", "-0.53***", "-0.51***", "0.12"),
    @("-0.43***", "-0.48***", "-0.44***", "-0.01"),
    @("-0.43***", "-0.48***", "-0.42***", "-0.1"),
    @("-0.31***", "-0.41***", "-0.32***", "-0.21*")
)
$dataRange = $ws.Range("B2:E6")
$dataRange.NumberFormat = "@"
for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $data[$r][$c]
    }
}
$dataRange.Style = "Normal"
